# Updated cryptos list: refresh prices and 1h volume changes,
# and correct the row order for Aptos/Litecoin and EthereumClassic/Fetch.AI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.319.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").Value = "'3.118.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'218.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.62%  "

$ws.Range("D6").Value = "'622.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("D7").Value = "'0.990"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +24.32%  "

$ws.Range("D8").Value = "'0.378"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.04%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'3.115.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("E11").Value = "  +20.21%  "

$ws.Range("D12").Value = "'0.191"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.38%  "

$ws.Range("D13").Value = "'0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.22%  "

$ws.Range("D14").Value = "'34.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.15%  "

$ws.Range("E15").Value = "  +2.64%  "

$ws.Range("D16").Value = "'91.096.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.66%  "

$ws.Range("D17").Value = "'3.682.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").Value = "'3.130.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("E19").Value = "  +13.53%  "

$ws.Range("D20").Value = "'0.0000219"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.11%  "

$ws.Range("D21").Value = "'14.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.74%  "

$ws.Range("D22").Value = "'435.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.81%  "

$ws.Range("D23").Value = "'8.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.79%  "

$ws.Range("D24").Value = "'5.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.68%  "

$ws.Range("D25").Value = "'6.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.99%  "

$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'12.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.56%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'86.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.03%  "

$ws.Range("D28").Value = "'3.274.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.48%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").Value = "'0.168"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.73%  "

$ws.Range("D31").Value = "'9.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.05%  "

$ws.Range("D32").Value = "'524.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.65%  "

$ws.Range("D33").Value = "'0.895"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.19%  "

$ws.Range("D34").Value = "'3.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.01%  "

$ws.Range("D35").Value = "'7.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.36%  "

$ws.Range("E36").Value = "  +13.16%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'23.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.87%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.78%  "

$ws.Range("E39").Value = "  +4.07%  "

$ws.Range("D40").Value = "'0.0896"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +32.36%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").Value = "'0.152"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.39%  "

$ws.Range("D44").Value = "'0.402"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.50%  "

$ws.Range("E46").Value = "  +6.94%  "

$ws.Range("D47").Value = "'148.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("D48").Value = "'44.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.74%  "

$ws.Range("E49").Value = "  +8.73%  "

$ws.Range("D50").Value = "'166.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.79%  "

$ws.Range("E51").Value = "  +7.89%  "
